$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the "Conversión del día" note with new rates ---
$ws1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 4.57 = 18142.03 pesos`n✅ 18142.03 pesos = 4.56 = 934.45 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$ws1.Range("A1").Value = $newText

# --- Sheet "tasas": update tasa values ---
$ws2 = $wb.Worksheets.Item("tasas")

$ws2.Range("N10").Value = 218.97
$ws2.Range("O10").Value = 3972.56
$ws2.Range("N12").Value = 3980
$ws2.Range("O12").Value = 205
